# Update the date/weekday heading at the top of the document.
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2025-05-09 Friday"

# Update the division-problem cells in the single worksheet table.
# Cells are addressed by (row, column) to avoid any ambiguity caused by
# values that coincide with other cells' old/new text (e.g. "37÷5=").
$t = $d.Tables.Item(1)

$updates = @(
    @(1, 1, "97÷5="),
    @(1, 2, "37÷5="),
    @(1, 3, "28÷8="),
    @(1, 4, "42÷6="),
    @(1, 5, "68÷5="),

    @(5, 1, "10÷6="),
    @(5, 2, "38÷7="),
    @(5, 3, "46÷9="),
    @(5, 4, "18÷6="),
    @(5, 5, "97÷8="),

    @(9, 1, "77÷7="),
    @(9, 2, "65÷2="),
    @(9, 3, "36÷2="),
    @(9, 4, "19÷5="),
    @(9, 5, "28÷4="),

    @(13, 1, "90÷8="),
    @(13, 2, "55÷3="),
    @(13, 3, "97÷4="),
    @(13, 4, "27÷6="),
    @(13, 5, "53÷9="),

    @(17, 1, "74÷4="),
    @(17, 2, "31÷4="),
    @(17, 3, "50÷2="),
    @(17, 4, "23÷5="),
    @(17, 5, "36÷3=")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $text = $u[2]
    $t.Cell($row, $col).Range.Text = $text
}
